$d = $word.ActiveDocument

# Locate an existing list paragraph that already uses the numId="3" bulleted
# list so the new bullet items below can continue that same list instance
# (rather than minting a brand-new numId).
$listTemplate = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Har spurgt mig selv hvordan patientens*") {
        $listTemplate = $para.Range.ListFormat.ListTemplate
        break
    }
}

# The document's trailing empty paragraph (just before the sectPr) becomes
# the new "11/12" Heading1 entry.
$pHeading = $d.Paragraphs.Last
$rHeading = $pHeading.Range
$rHeading.Text = "11/12 "
$rHeading.LanguageID = "da-DK"
$pHeading.Style = "Heading1"

# Subtitle paragraph.
$rHeading.InsertParagraphAfter()
$pSubtitle = $d.Paragraphs.Last
$rSubtitle = $pSubtitle.Range
$rSubtitle.Text = "Debugging af Activation Function"
$rSubtitle.LanguageID = "da-DK"
$pSubtitle.Style = "Subtitle"

# Bullet 1.
$rSubtitle.InsertParagraphAfter()
$pBullet1 = $d.Paragraphs.Last
$rBullet1 = $pBullet1.Range
$rBullet1.Text = "Activation Function fungerer i et test scenarie"
$rBullet1.LanguageID = "da-DK"
$pBullet1.Style = "ListParagraph"
$rBullet1.ListFormat.ApplyListTemplate($listTemplate, $true)

# Bullet 2.
$rBullet1.InsertParagraphAfter()
$pBullet2 = $d.Paragraphs.Last
$rBullet2 = $pBullet2.Range
$rBullet2.Text = "Den fungerer dog ikke i modellen"
$rBullet2.LanguageID = "da-DK"
$pBullet2.Style = "ListParagraph"
$rBullet2.ListFormat.ApplyListTemplate($listTemplate, $true)

# Bullet 3.
$rBullet2.InsertParagraphAfter()
$pBullet3 = $d.Paragraphs.Last
$rBullet3 = $pBullet3.Range
$rBullet3.Text = "Hypotese: Housekeeping er problemet"
$rBullet3.LanguageID = "da-DK"
$pBullet3.Style = "ListParagraph"
$rBullet3.ListFormat.ApplyListTemplate($listTemplate, $true)
